$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B ("Element"), shifting PC Number and
# everything after it one column to the right.
$ws.Columns("B:B").Insert()

$ws.Range("B1").Value = "Element"
$ws.Range("B2").Value = "MEP/N7101-1"
$ws.Range("B3").Value = "MEP/N7102-1"

$ws.Columns("B:B").ColumnWidth = 12.85546875
